$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CMP_Rep_Link")

# --- New validation columns H/I: list of ProjIDs (H) cross-checked against
#     column A via VLOOKUP (I). H12 = 24 is a ProjID that is not yet present
#     in A2:A38, so I12 legitimately resolves to #N/A.
$ws.Range("H6").Value = 2
$ws.Range("I6").Formula = '=VLOOKUP(H6,$A$2:$A$38,1,FALSE)'

$hValues = @(124,4,5,19,22,24,25,29,30,31,32,33,34,41,42,102,103,104,105,106,107,110,111,112,113,115,120,122)
$r = 7
foreach ($hv in $hValues) {
    $ws.Range("H" + $r).Value = $hv
    $r++
}
# Shared formula across I7:I34 (anchored relative to H7)
$ws.Range("I7:I34").Formula = '=VLOOKUP(H7,$A$2:$A$38,1,FALSE)'

# --- New row: the missing ProjID (24) gets added to the table, with a note
#     in column E explaining why it didn't resolve via VLOOKUP above.
$ws.Range("A39").Value = 24
$ws.Range("B39").Value = "Kapolei Interchange Complex Phase 2"
$ws.Range("C39").Value = 10987
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = "This is a new interchange, and can't accurately be captured by the regional model"

# --- Restore the cursor/selection that the author left on this sheet
#     without disturbing which sheet/tab is actually active in the workbook.
$ws.Activate()
$ws.Range("C24").Select()
$wb.Worksheets.Item("Proj Attributes and Scenarios").Activate()
